$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.710.11'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.22%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.598.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.44%  '

# Row 4
$ws.Range('E4').Value = '  +0.26%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.41%  '

# Row 6
$ws.Range('E6').Value = '  +0.44%  '

# Row 7
$ws.Range('E7').Value = '  +0.20%  '

# Row 8
$ws.Range('E8').Value = '  +0.66%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.247'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.00%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.49'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.34%  '

# Row 11
$ws.Range('E11').Value = '  +0.03%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.822.42'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.39%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.599.59'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.63%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.58%  '

# Row 15
$ws.Range('E15').Value = '  +0.79%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.18'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.20%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0₃0766'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.41%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.657.21'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.15%  '

# Row 19
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '209.58'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.09%  '

# Row 20
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.00'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.28%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.08'
$ws.Range('D21').Style = 'Normal'

# Row 22
$ws.Range('E22').Value = '  +1.38%  '

# Row 23
$ws.Range('E23').Value = '  +0.18%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.95'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.46%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.01'
$ws.Range('D25').Style = 'Normal'

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.18%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.12'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.57%  '

# Row 28
$ws.Range('E28').Value = '  +0.18%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.35'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.94%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0520'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.34%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.16'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.27%  '

# Row 32
$ws.Range('E32').Value = '  +0.85%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.97'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.54%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.285.51'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.52%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.617'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.26%  '

# Row 36
$ws.Range('E36').Value = '  -0.48%  '

# Row 37
$ws.Range('E37').Value = '  +0.38%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0171'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.13%  '

# Row 39
$ws.Range('E39').Value = '  +17.90%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.827'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.05%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.44'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.96%  '

# Row 42
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.785'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.01%  '

# Row 43
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.19'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.52%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.24'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.19%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.734.40'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.34%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.43%  '

# Row 47
$ws.Range('E47').Value = '  -2.12%  '

# Row 48
$ws.Range('E48').Value = '  +0.54%  '

# Row 49
$ws.Range('E49').Value = '  +0.70%  '

# Row 50
$ws.Range('E50').Value = '  +0.20%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.32'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.31%  '
